# Auto-generated Excel COM-interop script to apply the Anima_Profits.xlsx diff
# across the ALC, ARM, BSM, CRP, CUL, LTW, WVR worksheets.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1074.1
$ws.Range("I2").Value = 1174.5555
$ws.Range("J2").Value = 170
$ws.Range("K2").Value = 1174.5555
$ws.Range("L2").Value = 170
$ws.Range("M2").Value = -1061.5555
$ws.Range("N2").Value = -396
$ws.Range("H92").Value = 11112688
$ws.Range("I92").Value = 17544876
$ws.Range("J92").Value = 2545.3635
$ws.Range("K92").Value = 17544876
$ws.Range("L92").Value = 2545.3635
$ws.Range("M92").Value = -17543628
$ws.Range("N92").Value = -5041.363499999999
$ws.Range("H103").Value = 53799.26
$ws.Range("I103").Value = 100519.8
$ws.Range("J103").Value = 1887.5555
$ws.Range("K103").Value = 301559.4
$ws.Range("L103").Value = 5662.666499999999
$ws.Range("M103").Value = -300973.4
$ws.Range("N103").Value = -6834.666499999999
$ws.Range("H129").Value = 1172.7261
$ws.Range("I129").Value = 684.5714
$ws.Range("J129").Value = 1224.5
$ws.Range("K129").Value = 2053.7142
$ws.Range("L129").Value = 3673.5
$ws.Range("M129").Value = 2946.2858
$ws.Range("N129").Value = -13673.5
$ws.Range("H132").Value = 2976.5625
$ws.Range("I132").Value = 2975.8064
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8927.4192
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6397.4192
$ws.Range("N132").Value = -14060
$ws.Range("H134").Value = 79726.664
$ws.Range("J134").Value = 79726.664
$ws.Range("L134").Value = 79726.664
$ws.Range("N134").Value = -89866.664
$ws.Range("H138").Value = 3676.1208
$ws.Range("J138").Value = 3878.75
$ws.Range("L138").Value = 11636.25
$ws.Range("N138").Value = -21916.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 19186.5
$ws.Range("J80").Value = 19186.5
$ws.Range("L80").Value = 19186.5
$ws.Range("N80").Value = -21182.5
$ws.Range("H83").Value = 19186.5
$ws.Range("J83").Value = 19186.5
$ws.Range("L83").Value = 57559.5
$ws.Range("N83").Value = -67543.5
$ws.Range("H93").Value = 57500
$ws.Range("J93").Value = 57500
$ws.Range("L93").Value = 57500
$ws.Range("N93").Value = -62492
$ws.Range("H107").Value = 51409
$ws.Range("J107").Value = 51409
$ws.Range("L107").Value = 51409
$ws.Range("N107").Value = -59089
$ws.Range("H108").Value = 59936.8
$ws.Range("J108").Value = 59936.8
$ws.Range("L108").Value = 59936.8
$ws.Range("N108").Value = -67616.8
$ws.Range("H115").Value = 77500
$ws.Range("J115").Value = 77500
$ws.Range("L115").Value = 77500
$ws.Range("N115").Value = -80634
$ws.Range("H132").Value = 3320.02
$ws.Range("I132").Value = 2417.7778
$ws.Range("J132").Value = 5640.0713
$ws.Range("K132").Value = 7253.3334
$ws.Range("L132").Value = 16920.2139
$ws.Range("M132").Value = -4723.3334
$ws.Range("N132").Value = -21980.2139

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 22926.666
$ws.Range("J50").Value = 22926.666
$ws.Range("L50").Value = 22926.666
$ws.Range("N50").Value = -24074.666
$ws.Range("H93").Value = 62224
$ws.Range("J93").Value = 62224
$ws.Range("L93").Value = 62224
$ws.Range("N93").Value = -65968
$ws.Range("H96").Value = 35429.816
$ws.Range("I96").Value = 8714
$ws.Range("J96").Value = 41366.668
$ws.Range("K96").Value = 8714
$ws.Range("L96").Value = 41366.668
$ws.Range("M96").Value = -5968
$ws.Range("N96").Value = -46858.668
$ws.Range("H97").Value = 16446.273
$ws.Range("I97").Value = 4937.5557
$ws.Range("K97").Value = 4937.5557
$ws.Range("M97").Value = -3946.5557
$ws.Range("H109").Value = 48335.5
$ws.Range("J109").Value = 48335.5
$ws.Range("L109").Value = 48335.5
$ws.Range("N109").Value = -51109.5
$ws.Range("H115").Value = 78040.30499999999
$ws.Range("J115").Value = 78040.30499999999
$ws.Range("L115").Value = 78040.30499999999
$ws.Range("N115").Value = -81174.30499999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 419.04
$ws.Range("I22").Value = 296
$ws.Range("J22").Value = 532.61536
$ws.Range("K22").Value = 296
$ws.Range("L22").Value = 532.61536
$ws.Range("M22").Value = 54
$ws.Range("N22").Value = -1232.61536
$ws.Range("H31").Value = 5598.26
$ws.Range("I31").Value = 1452.4584
$ws.Range("J31").Value = 7475.604
$ws.Range("K31").Value = 1452.4584
$ws.Range("L31").Value = 7475.604
$ws.Range("M31").Value = -1157.4584
$ws.Range("N31").Value = -8065.604
$ws.Range("H34").Value = 5598.26
$ws.Range("I34").Value = 1452.4584
$ws.Range("J34").Value = 7475.604
$ws.Range("K34").Value = 1452.4584
$ws.Range("L34").Value = 7475.604
$ws.Range("M34").Value = -1250.4584
$ws.Range("N34").Value = -7879.604
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null
$ws.Range("H127").Value = 51890
$ws.Range("J127").Value = 51890
$ws.Range("L127").Value = 51890
$ws.Range("N127").Value = -61810
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H140").Value = 64770
$ws.Range("J140").Value = 64770
$ws.Range("L140").Value = 64770
$ws.Range("N140").Value = -75130

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 798722
$ws.Range("I129").Value = 451
$ws.Range("J129").Value = 1167154.8
$ws.Range("K129").Value = 1353
$ws.Range("L129").Value = 3501464.4
$ws.Range("M129").Value = 3647
$ws.Range("N129").Value = -3511464.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 29325
$ws.Range("J63").Value = 29325
$ws.Range("L63").Value = 29325
$ws.Range("N63").Value = -30823
$ws.Range("H66").Value = 29325
$ws.Range("J66").Value = 29325
$ws.Range("L66").Value = 87975
$ws.Range("N66").Value = -95463
$ws.Range("H68").Value = 2666.6667
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 3250
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 3250
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -4748
$ws.Range("H70").Value = 25290.75
$ws.Range("J70").Value = 25290.75
$ws.Range("L70").Value = 25290.75
$ws.Range("N70").Value = -25830.75
$ws.Range("H71").Value = 2666.6667
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 3250
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 16250
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -23738
$ws.Range("H73").Value = 25290.75
$ws.Range("J73").Value = 25290.75
$ws.Range("L73").Value = 25290.75
$ws.Range("N73").Value = -27162.75
$ws.Range("H112").Value = 69990
$ws.Range("J112").Value = 69990
$ws.Range("L112").Value = 69990
$ws.Range("N112").Value = -72944
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = $null
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820
$ws.Range("H130").Value = 69000
$ws.Range("J130").Value = 69000
$ws.Range("L130").Value = 69000
$ws.Range("N130").Value = -79040

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 58000
$ws.Range("J27").Value = 58000
$ws.Range("L27").Value = 58000
$ws.Range("N27").Value = -58138
$ws.Range("H68").Value = 95000
$ws.Range("J68").Value = 95000
$ws.Range("L68").Value = 95000
$ws.Range("N68").Value = -96622
$ws.Range("H71").Value = 95000
$ws.Range("J71").Value = 95000
$ws.Range("L71").Value = 285000
$ws.Range("N71").Value = -293112
$ws.Range("H80").Value = 60720.2
$ws.Range("J80").Value = 60720.2
$ws.Range("L80").Value = 60720.2
$ws.Range("N80").Value = -62716.2
$ws.Range("H83").Value = 60720.2
$ws.Range("J83").Value = 60720.2
$ws.Range("L83").Value = 182160.6
$ws.Range("N83").Value = -192144.6
$ws.Range("H102").Value = 42800
$ws.Range("J102").Value = 42800
$ws.Range("L102").Value = 42800
$ws.Range("N102").Value = -49290
$ws.Range("H103").Value = 56800
$ws.Range("J103").Value = 56800
$ws.Range("L103").Value = 56800
$ws.Range("N103").Value = -59144
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = $null
